$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '276.15'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '6.264'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06197'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.560'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.533'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '6.579'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8259'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1667'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08316'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03513'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03183'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09161'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.766'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001628'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04690'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.006280'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.006218'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.001068'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0001501'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.322'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.01398'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.3291'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1242'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0002737'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04748'
$ws.Range("B41").Value = 'CEJI'
$ws.Range("C41").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005302'
$ws.Range("E41").Value = '40CEJICEJI'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.007066'
$ws.Range("E42").Value = '41KickTokenKICK'
$ws.Range("B43").Value = 'BKEXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1119'
$ws.Range("E43").Value = '42BKEXTokenBKK'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01135'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006384'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.7232'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.001400'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00001901'
